$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at O:P, shifting the old O:U (and their data)
# two columns to the right, turning them into Q:W. This also grows the
# sheet's dimension from A1:U6 to A1:W6 and keeps the header style/border
# (s="1") applied to the newly inserted header cells because they inherit
# formatting from the columns immediately to their left.
$ws.Range("O1:P1").EntireColumn.Insert()

# Rename the two predicate-detection headers that keep their position.
$ws.Range("M1").Value = "Detected Predicates Doc Parent"
$ws.Range("N1").Value = "Detected Predicates Doc Related"

# Fill in the headers for the two brand-new columns.
$ws.Range("O1").Value = "Correct Pred Predicates Parents"
$ws.Range("P1").Value = "Correct Pred Predicates Related"

# Populate the new "Correct Pred Predicates Parents/Related" data columns.
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1

$ws.Range("O3").Value = 4
$ws.Range("P3").Value = 4

$ws.Range("O4").Value = 2
$ws.Range("P4").Value = 2

$ws.Range("O5").Value = 4
$ws.Range("P5").Value = 4

$ws.Range("O6").Value = 2
$ws.Range("P6").Value = 2
